$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1891.0834
$ws.Range("I92").Value = 1881.1818
$ws.Range("K92").Value = 1881.1818
$ws.Range("M92").Value = -633.1818000000001

$ws.Range("H96").Value = 2333.3333
$ws.Range("I96").Value = 2333.3333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 6999.999899999999
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -5626.999899999999
$ws.Range("N96").ClearContents()

$ws.Range("H106").Value = 2696.9524
$ws.Range("I106").Value = 3268.2307
$ws.Range("J106").Value = 1768.625
$ws.Range("K106").Value = 3268.2307
$ws.Range("L106").Value = 1768.625
$ws.Range("M106").Value = -2637.2307
$ws.Range("N106").Value = -3030.625

$ws.Range("H113").Value = 13601.5
$ws.Range("I113").Value = 6048.1665
$ws.Range("K113").Value = 6048.1665
$ws.Range("M113").Value = -2794.1665

$ws.Range("H115").Value = 791.6667
$ws.Range("I115").Value = 791.6667
$ws.Range("K115").Value = 2375.0001
$ws.Range("M115").Value = -808.0001000000002

$ws.Range("H116").Value = 4453807
$ws.Range("I116").Value = 6544699
$ws.Range("J116").Value = 10661.875
$ws.Range("K116").Value = 6544699
$ws.Range("L116").Value = 10661.875
$ws.Range("M116").Value = -6541257
$ws.Range("N116").Value = -17545.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1853212.5
$ws.Range("I32").Value = 624.0599999999999
$ws.Range("J32").Value = 11116154
$ws.Range("K32").Value = 624.0599999999999
$ws.Range("L32").Value = 11116154
$ws.Range("M32").Value = -337.0599999999999
$ws.Range("N32").Value = -11116728

$ws.Range("H45").Value = 1491.5714
$ws.Range("I45").Value = 1554.625
$ws.Range("K45").Value = 1554.625
$ws.Range("M45").Value = -1177.625

$ws.Range("H61").Value = 2767.12
$ws.Range("I61").Value = 1847.8
$ws.Range("K61").Value = 1847.8
$ws.Range("M61").Value = -1635.8

$ws.Range("H74").Value = 1461.5
$ws.Range("I74").Value = 1363.2142
$ws.Range("J74").Value = 2149.5
$ws.Range("K74").Value = 1363.2142
$ws.Range("L74").Value = 2149.5
$ws.Range("M74").Value = -489.2141999999999
$ws.Range("N74").Value = -3897.5

$ws.Range("H77").Value = 1461.5
$ws.Range("I77").Value = 1363.2142
$ws.Range("J77").Value = 2149.5
$ws.Range("K77").Value = 6816.071
$ws.Range("L77").Value = 10747.5
$ws.Range("M77").Value = -2448.071
$ws.Range("N77").Value = -19483.5

$ws.Range("H88").Value = 5591.6
$ws.Range("I88").Value = 1878.1428
$ws.Range("J88").Value = 8840.875
$ws.Range("K88").Value = 1878.1428
$ws.Range("L88").Value = 8840.875
$ws.Range("M88").Value = -1472.1428
$ws.Range("N88").Value = -9652.875

$ws.Range("H91").Value = 5591.6
$ws.Range("I91").Value = 1878.1428
$ws.Range("J91").Value = 8840.875
$ws.Range("K91").Value = 1878.1428
$ws.Range("L91").Value = 8840.875
$ws.Range("M91").Value = -474.1428000000001
$ws.Range("N91").Value = -11648.875

$ws.Range("H101").Value = 32499.5
$ws.Range("J101").Value = 32499.5
$ws.Range("L101").Value = 32499.5
$ws.Range("N101").Value = -38989.5

$ws.Range("H122").Value = 2194.4583
$ws.Range("I122").Value = 1142.8667
$ws.Range("K122").Value = 3428.6001
$ws.Range("M122").Value = -978.6001000000001

$ws.Range("H132").Value = 978710.1
$ws.Range("I132").Value = 1265072.2
$ws.Range("K132").Value = 3795216.6
$ws.Range("M132").Value = -3792686.6

$ws.Range("H136").Value = 2767.12
$ws.Range("I136").Value = 1847.8
$ws.Range("K136").Value = 5543.4
$ws.Range("M136").Value = -2993.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8931089
$ws.Range("I20").Value = 20411104
$ws.Range("K20").Value = 20411104
$ws.Range("M20").Value = -20410857

$ws.Range("H82").Value = 32460.285
$ws.Range("J82").Value = 26805.5
$ws.Range("L82").Value = 26805.5
$ws.Range("N82").Value = -27571.5

$ws.Range("H85").Value = 32460.285
$ws.Range("J85").Value = 26805.5
$ws.Range("L85").Value = 26805.5
$ws.Range("N85").Value = -29457.5

$ws.Range("H97").Value = 18705.428
$ws.Range("J97").Value = 23647
$ws.Range("L97").Value = 23647
$ws.Range("N97").Value = -25629

$ws.Range("H105").Value = 4531.8
$ws.Range("J105").Value = 1447.5
$ws.Range("L105").Value = 1447.5
$ws.Range("N105").Value = -4941.5

$ws.Range("H134").Value = 1838112.2
$ws.Range("I134").Value = 2167770.2
$ws.Range("J134").Value = 24993.25
$ws.Range("K134").Value = 6503310.600000001
$ws.Range("L134").Value = 74979.75
$ws.Range("M134").Value = -6500775.600000001
$ws.Range("N134").Value = -80049.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 500
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 1500
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -2076

$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2831

$ws.Range("H37").Value = 218564
$ws.Range("J37").Value = 218564
$ws.Range("L37").Value = 655692
$ws.Range("N37").Value = -655916

$ws.Range("H61").Value = 757.0909
$ws.Range("J61").Value = 1189.8
$ws.Range("L61").Value = 3569.4
$ws.Range("N61").Value = -3999.4

$ws.Range("H70").Value = 5402
$ws.Range("I70").Value = 3502.75
$ws.Range("K70").Value = 10508.25
$ws.Range("M70").Value = -10193.25

$ws.Range("H73").Value = 5402
$ws.Range("I73").Value = 3502.75
$ws.Range("K73").Value = 10508.25
$ws.Range("M73").Value = -9416.25

$ws.Range("H123").Value = 29500
$ws.Range("J123").Value = 29500
$ws.Range("L123").Value = 88500
$ws.Range("N123").Value = -93400

$ws.Range("H131").Value = 37683384
$ws.Range("J131").Value = 20836146
$ws.Range("L131").Value = 62508438
$ws.Range("N131").Value = -62518518

$ws.Range("H137").Value = 4549.1055
$ws.Range("I137").Value = 5552.4546
$ws.Range("J137").Value = 3169.5
$ws.Range("K137").Value = 16657.3638
$ws.Range("L137").Value = 9508.5
$ws.Range("M137").Value = -11557.3638
$ws.Range("N137").Value = -19708.5

$ws.Range("H139").Value = 2313.3076
$ws.Range("I139").Value = 2006.6364
$ws.Range("K139").Value = 6019.9092
$ws.Range("M139").Value = -879.9092000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4245
$ws.Range("I102").Value = 3756.0908
$ws.Range("K102").Value = 3756.0908
$ws.Range("M102").Value = -2134.0908

$ws.Range("H122").Value = 8814.92
$ws.Range("I122").Value = 8506.637000000001
$ws.Range("K122").Value = 25519.911
$ws.Range("M122").Value = -23069.911

$ws.Range("H132").Value = 62504640
$ws.Range("I132").Value = 111114990
$ws.Range("K132").Value = 333344970
$ws.Range("M132").Value = -333342440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 589.6
$ws.Range("J22").Value = 583
$ws.Range("L22").Value = 583
$ws.Range("N22").Value = -1173

$ws.Range("H27").Value = 589.6
$ws.Range("J27").Value = 583
$ws.Range("L27").Value = 583
$ws.Range("N27").Value = -797

$ws.Range("H40").Value = 7401
$ws.Range("I40").Value = 6163.6665
$ws.Range("K40").Value = 6163.6665
$ws.Range("M40").Value = -6027.6665

$ws.Range("H55").Value = 2320.1365
$ws.Range("I55").Value = 629.2
$ws.Range("J55").Value = 3729.25
$ws.Range("K55").Value = 629.2
$ws.Range("L55").Value = 3729.25
$ws.Range("M55").Value = -456.2
$ws.Range("N55").Value = -4075.25
